# Rename existing "Sheet1" to "Variables"
$wb = $excel.ActiveWorkbook
$wsVariables = $wb.Worksheets.Item("Sheet1")
$wsVariables.Name = "Variables"

# Add a new worksheet "removeNaNs" after "Variables"
$wsRemoveNaNs = $wb.Worksheets.Add($null, $wsVariables)
$wsRemoveNaNs.Name = "removeNaNs"

# Populate the new sheet with header + values
$wsRemoveNaNs.Range("A1").Value = "Variable"
$wsRemoveNaNs.Range("A1").Font.Bold = $true
$wsRemoveNaNs.Range("A2").Value = "Ctr"
$wsRemoveNaNs.Range("A3").Value = "PT_100"

# Add autofilter on the new sheet
$wsRemoveNaNs.Range("A1").AutoFilter() | Out-Null

# Set selection/view state to match target
$wsVariables.Range("A7").Select()
$wsRemoveNaNs.Range("B5").Select()

# Activate the removeNaNs sheet (tab selected)
$wsRemoveNaNs.Activate()
